function Set-TextValue($cell, $value) {
    # Force the cell to keep a numeric-looking string (e.g. "011765" / "1.01")
    # as literal text instead of letting Excel auto-convert it to a number
    # (and losing the leading zero / trailing formatting in the process).
    # ClearFormats() afterwards drops the temporary "@" text format so the
    # cell is left with no explicit style, matching the plain data cells
    # elsewhere in the workbook.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Before: sheets are [2021-Q3, 2021-Q4, 总计]            (总计 = rId3/sheetId3)
# After : sheets are [2021-Q3, 2021-Q4, 2022-Q1, 总计]    (2022-Q1 keeps
#          sheetId3/rId3 and gets new fund-level detail rows; a brand-new
#          总计 sheet (sheetId4/rId4) is appended with the running
#          quarter-by-quarter roll-up, now including 2022-Q1).
# ---------------------------------------------------------------------------

# Step 1: duplicate the current "总计" sheet (index 3). Excel places the
# copy right after the original, so we end up with:
#   [2021-Q3, 2021-Q4, 总计, 总计 (2)]
$totalSheet = $wb.Worksheets.Item(3)
$totalSheet.Copy([System.Reflection.Missing]::Value, $totalSheet)

$q1Sheet = $wb.Worksheets.Item(3)
$newTotalSheet = $wb.Worksheets.Item(4)

# Rename the original first (otherwise renaming the copy collides with the
# still-existing "总计" name).
$q1Sheet.Name = "2022-Q1"
$newTotalSheet.Name = "总计"

# ---------------------------------------------------------------------------
# Step 2: rewrite "2022-Q1" (previously the "总计" sheet) with per-fund
# holding details, matching the layout used by 2021-Q3 / 2021-Q4.
# ---------------------------------------------------------------------------

# Clear out the old 总计-style content (B1:D3) before laying out the new grid.
$q1Sheet.Range("A1:H3").ClearContents()

$q1Sheet.Cells.Item(1, 2).Value = "基金代码"
$q1Sheet.Cells.Item(1, 3).Value = "基金名称"
$q1Sheet.Cells.Item(1, 4).Value = "基金规模"
$q1Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q1Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q1Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1Sheet.Cells.Item(1, 8).Value = "仓位排名"

# B1:D1 already carry the bold/bordered/centered header style inherited from
# the old sheet; copy that same formatting across the newly added header
# cells E1:H1.
$q1Sheet.Cells.Item(1, 2).Copy()
$q1Sheet.Range($q1Sheet.Cells.Item(1, 5), $q1Sheet.Cells.Item(1, 8)).PasteSpecial(-4122)

# Row 2 -- 011765 兴银高端制造混合A
$q1Sheet.Cells.Item(2, 1).Value = 0
Set-TextValue $q1Sheet.Cells.Item(2, 2) "011765"
$q1Sheet.Cells.Item(2, 3).Value = "兴银高端制造混合A"
Set-TextValue $q1Sheet.Cells.Item(2, 4) "1.01"
Set-TextValue $q1Sheet.Cells.Item(2, 5) "93.23"
Set-TextValue $q1Sheet.Cells.Item(2, 6) "2.60"
Set-TextValue $q1Sheet.Cells.Item(2, 7) "0.0263"
$q1Sheet.Cells.Item(2, 8).Value = 8

# Row 3 -- 011766 兴银高端制造混合C
$q1Sheet.Cells.Item(3, 1).Value = 1
Set-TextValue $q1Sheet.Cells.Item(3, 2) "011766"
$q1Sheet.Cells.Item(3, 3).Value = "兴银高端制造混合C"
Set-TextValue $q1Sheet.Cells.Item(3, 4) "0.39"
Set-TextValue $q1Sheet.Cells.Item(3, 5) "93.23"
Set-TextValue $q1Sheet.Cells.Item(3, 6) "2.60"
Set-TextValue $q1Sheet.Cells.Item(3, 7) "0.0101"
$q1Sheet.Cells.Item(3, 8).Value = 8

# A2/A3 (the 0 / 1 row-index column) keep the bold/bordered/centered look too.
$q1Sheet.Cells.Item(2, 1).Copy()
$q1Sheet.Cells.Item(3, 1).PasteSpecial(-4122)
$q1Sheet.Cells.Item(2, 1).Value = 0
$q1Sheet.Cells.Item(3, 1).Value = 1

# ---------------------------------------------------------------------------
# Step 3: update the new "总计" sheet -- push the existing rows down one and
# insert the new 2022-Q1 roll-up row at the top of the data.
# ---------------------------------------------------------------------------
$newTotalSheet.Cells.Item(4, 1).Value = 2
$newTotalSheet.Cells.Item(4, 2).Value = "2021-Q3"
$newTotalSheet.Cells.Item(4, 3).Value = 6
$newTotalSheet.Cells.Item(4, 4).Value = 0.25

$newTotalSheet.Cells.Item(3, 1).Value = 1
$newTotalSheet.Cells.Item(3, 2).Value = "2021-Q4"
$newTotalSheet.Cells.Item(3, 3).Value = 5
$newTotalSheet.Cells.Item(3, 4).Value = 0.61

$newTotalSheet.Cells.Item(2, 1).Value = 0
$newTotalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$newTotalSheet.Cells.Item(2, 3).Value = 2
$newTotalSheet.Cells.Item(2, 4).Value = 0.04

# Give the newly-created A4 the same bold/bordered/centered style as A2/A3.
$newTotalSheet.Cells.Item(3, 1).Copy()
$newTotalSheet.Cells.Item(4, 1).PasteSpecial(-4122)
$newTotalSheet.Cells.Item(4, 1).Value = 2
